$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 51 (hunk 0)
$ws.Range("H51").Value = 1785.7142
$ws.Range("I51").Value = 1775
$ws.Range("J51").Value = 1800
$ws.Range("K51").Value = 1775
$ws.Range("L51").Value = 1800
$ws.Range("M51").Value = -1291
$ws.Range("N51").Value = -2768
# Row 105 (hunk 1)
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
# Row 107 (hunk 2)
$ws.Range("H107").Value = 1501.1111
$ws.Range("I107").Value = 1563.75
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1563.75
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 356.25
$ws.Range("N107").Value = -4840
# Row 137 (hunk 3)
$ws.Range("H137").Value = 2875.8286
$ws.Range("I137").Value = 2719.818
$ws.Range("J137").Value = 5450
$ws.Range("K137").Value = 8159.454000000001
$ws.Range("L137").Value = 16350
$ws.Range("M137").Value = -5609.454000000001
$ws.Range("N137").Value = -21450
# Row 141 (hunk 4)
$ws.Range("H141").Value = 495843.38
$ws.Range("I141").Value = 1593.8235
$ws.Range("J141").Value = 1546123.6
$ws.Range("K141").Value = 4781.470499999999
$ws.Range("L141").Value = 4638370.800000001
$ws.Range("M141").Value = 398.5295000000006
$ws.Range("N141").Value = -4648730.800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45 (hunk 5)
$ws.Range("H45").Value = 1347.3513
$ws.Range("I45").Value = 1079.0344
$ws.Range("J45").Value = 2320
$ws.Range("K45").Value = 1079.0344
$ws.Range("L45").Value = 2320
$ws.Range("M45").Value = -702.0344
$ws.Range("N45").Value = -3074
# Row 61 (hunk 6)
$ws.Range("H61").Value = 2108.75
$ws.Range("I61").Value = 1230.5
$ws.Range("J61").Value = 6500
$ws.Range("K61").Value = 1230.5
$ws.Range("L61").Value = 6500
$ws.Range("M61").Value = -1018.5
$ws.Range("N61").Value = -6924
# Row 122 (hunk 7)
$ws.Range("H122").Value = 2722.7058
$ws.Range("I122").Value = 1909.2
$ws.Range("J122").Value = 3884.8572
$ws.Range("K122").Value = 5727.6
$ws.Range("L122").Value = 11654.5716
$ws.Range("M122").Value = -3277.6
$ws.Range("N122").Value = -16554.5716
# Row 132 (hunk 8)
$ws.Range("H132").Value = 16131376
$ws.Range("I132").Value = 24391536
$ws.Range("J132").Value = 4394.857
$ws.Range("K132").Value = 73174608
$ws.Range("L132").Value = 13184.571
$ws.Range("M132").Value = -73172078
$ws.Range("N132").Value = -18244.571
# Row 136 (hunk 9)
$ws.Range("H136").Value = 2108.75
$ws.Range("I136").Value = 1230.5
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 3691.5
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -1141.5
$ws.Range("N136").Value = -24600
# Row 141 (hunk 10)
$ws.Range("H141").Value = 29496.25
$ws.Range("J141").Value = 29496.25
$ws.Range("L141").Value = 29496.25
$ws.Range("N141").Value = -39856.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 5 (hunk 11)
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
# Row 134 (hunk 12)
$ws.Range("H134").Value = 2252.88
$ws.Range("I134").Value = 1414.7
$ws.Range("K134").Value = 4244.1
$ws.Range("M134").Value = -1709.1

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (hunk 13)
$ws.Range("H22").Value = 3675
$ws.Range("J22").Value = 4666.6665
$ws.Range("L22").Value = 4666.6665
$ws.Range("N22").Value = -5366.6665
# Row 75 (hunk 14)
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78 (hunk 15)
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 132 (hunk 16)
$ws.Range("H132").Value = 2021.1875
$ws.Range("I132").Value = 1482.925
$ws.Range("J132").Value = 4712.5
$ws.Range("K132").Value = 4448.775
$ws.Range("L132").Value = 14137.5
$ws.Range("M132").Value = -1918.775
$ws.Range("N132").Value = -19197.5
# Row 134 (hunk 17)
$ws.Range("H134").Value = 1610
$ws.Range("I134").Value = 871.5
$ws.Range("J134").Value = 2533.125
$ws.Range("K134").Value = 2614.5
$ws.Range("L134").Value = 7599.375
$ws.Range("M134").Value = -79.5
$ws.Range("N134").Value = -12669.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 59 (hunk 18)
$ws.Range("H59").Value = 1589.1666
$ws.Range("I59").Value = 211.66667
$ws.Range("J59").Value = 2966.6667
$ws.Range("K59").Value = 635.00001
$ws.Range("L59").Value = 8900.000100000001
$ws.Range("M59").Value = -95.00000999999997
$ws.Range("N59").Value = -9980.000100000001
# Row 64 (hunk 19)
$ws.Range("H64").Value = 1918.4
$ws.Range("I64").Value = 982.2857
$ws.Range("J64").Value = 2737.5
$ws.Range("K64").Value = 2946.8571
$ws.Range("L64").Value = 8212.5
$ws.Range("M64").Value = -2676.8571
$ws.Range("N64").Value = -8752.5
# Row 67 (hunk 20)
$ws.Range("H67").Value = 1918.4
$ws.Range("I67").Value = 982.2857
$ws.Range("J67").Value = 2737.5
$ws.Range("K67").Value = 2946.8571
$ws.Range("L67").Value = 8212.5
$ws.Range("M67").Value = -2010.8571
$ws.Range("N67").Value = -10084.5
# Row 116 (hunk 21)
$ws.Range("H116").Value = 2466.6667
$ws.Range("I116").Value = 1700
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 5100
$ws.Range("L116").Value = 12000
$ws.Range("M116").Value = -1658
$ws.Range("N116").Value = -18884
# Row 131 (hunk 22)
$ws.Range("H131").Value = 1115.3269
$ws.Range("J131").Value = 1005.5349
$ws.Range("L131").Value = 3016.6047
$ws.Range("N131").Value = -13096.6047
# Row 133 (hunk 23)
$ws.Range("H133").Value = 7289.6
$ws.Range("I133").Value = 7289.6
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 21868.8
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -16808.8
$ws.Range("N133").ClearContents()
# Row 134 (hunk 24)
$ws.Range("H134").Value = 2575.85
$ws.Range("I134").Value = 1235.2222
$ws.Range("J134").Value = 3672.7273
$ws.Range("K134").Value = 3705.6666
$ws.Range("L134").Value = 11018.1819
$ws.Range("M134").Value = 1364.3334
$ws.Range("N134").Value = -21158.1819

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 97 (hunk 25)
$ws.Range("H97").Value = 3522
$ws.Range("I97").Value = 2036.6666
$ws.Range("J97").Value = 5750
$ws.Range("K97").Value = 2036.6666
$ws.Range("L97").Value = 5750
$ws.Range("M97").Value = -1540.6666
$ws.Range("N97").Value = -6742
# Row 102 (hunk 26)
$ws.Range("H102").Value = 86845.336
$ws.Range("I102").Value = 2905
$ws.Range("J102").Value = 338666.34
$ws.Range("K102").Value = 2905
$ws.Range("L102").Value = 338666.34
$ws.Range("M102").Value = -1283
$ws.Range("N102").Value = -341910.34
# Row 122 (hunk 27)
$ws.Range("H122").Value = 4654.227
$ws.Range("I122").Value = 2278.125
$ws.Range("J122").Value = 6012
$ws.Range("K122").Value = 6834.375
$ws.Range("L122").Value = 18036
$ws.Range("M122").Value = -4384.375
$ws.Range("N122").Value = -22936
# Row 132 (hunk 28)
$ws.Range("H132").Value = 3213.7073
$ws.Range("I132").Value = 2775.5386
$ws.Range("J132").Value = 3973.2
$ws.Range("K132").Value = 8326.6158
$ws.Range("L132").Value = 11919.6
$ws.Range("M132").Value = -5796.6158
$ws.Range("N132").Value = -16979.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 29)
$ws.Range("H7").Value = 2600
$ws.Range("I7").Value = 1900
$ws.Range("K7").Value = 1900
$ws.Range("M7").Value = -1788
# Row 40 (hunk 30)
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2272
$ws.Range("M40").ClearContents()
# Row 122 (hunk 31)
$ws.Range("H122").Value = 3825
$ws.Range("I122").Value = 2728.5715
$ws.Range("J122").Value = 4677.778
$ws.Range("K122").Value = 8185.7145
$ws.Range("L122").Value = 14033.334
$ws.Range("M122").Value = -5735.7145
$ws.Range("N122").Value = -18933.334
# Row 126 (hunk 32)
$ws.Range("H126").Value = 2600
$ws.Range("I126").Value = 1900
$ws.Range("K126").Value = 5700
$ws.Range("M126").Value = -3230

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 70 (hunk 33)
$ws.Range("H70").Value = 90105
$ws.Range("J70").Value = 90105
$ws.Range("L70").Value = 90105
$ws.Range("N70").Value = -90735
# Row 73 (hunk 34)
$ws.Range("H73").Value = 90105
$ws.Range("J73").Value = 90105
$ws.Range("L73").Value = 90105
$ws.Range("N73").Value = -92289
# Row 81 (hunk 35)
$ws.Range("H81").Value = 907.6667
$ws.Range("I81").Value = 907.6667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1815.3334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -754.3334
$ws.Range("N81").ClearContents()
# Row 84 (hunk 36)
$ws.Range("H84").Value = 907.6667
$ws.Range("I84").Value = 907.6667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9076.666999999999
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3772.666999999999
$ws.Range("N84").ClearContents()
# Row 113 (hunk 37)
$ws.Range("H113").Value = 1364.9445
$ws.Range("I113").Value = 477.27274
$ws.Range("J113").Value = 2759.8572
$ws.Range("K113").Value = 1431.81822
$ws.Range("L113").Value = 8279.571599999999
$ws.Range("M113").Value = 738.1817799999999
$ws.Range("N113").Value = -12619.5716
# Row 122 (hunk 38)
$ws.Range("H122").Value = 717237.4399999999
$ws.Range("I122").Value = 1430675.1
$ws.Range("J122").Value = 3799.7144
$ws.Range("K122").Value = 4292025.300000001
$ws.Range("L122").Value = 11399.1432
$ws.Range("M122").Value = -4289575.300000001
$ws.Range("N122").Value = -16299.1432
# Row 140 (hunk 39)
$ws.Range("H140").Value = 38444.082
$ws.Range("J140").Value = 38444.082
$ws.Range("L140").Value = 38444.082
$ws.Range("N140").Value = -48804.082
